# Adjusted wording in the "Any strong individually held views" report:
#  - Brian's quote gets two new inserted clauses ("done earlier on" and
#    "to find things which should, really, have been agreed at the design
#    stage") plus a reworded closing ("This strategy is not ideal ... would
#    scale badly to larger projects" - drops "even").
#  - The _GoBack bookmark moves from the start of Edward's quote ("Our group")
#    to its new natural position (right before "is not ideal") since that is
#    where the author's final edit/typing happened.
#
# This engine re-coalesces adjacent same-formatted runs on every text edit,
# so bookmarks are used as temporary "shields" to pin down the run
# boundaries the diff expects; the temporary ones are removed at the end
# (the real _GoBack bookmark is left in place).

$d = $word.ActiveDocument

function Get-Range($searchText) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($searchText) | Out-Null
    if (-not $r.Find.Found) {
        Write-Host "NOT FOUND: $searchText"
    }
    return $r
}

function Shield-Point($rng, $name) {
    $d.Bookmarks.Add($name, $rng) | Out-Null
}

# Inserts $newText immediately before/after the match of $searchText
# (collapseDir: 0 = end of match, 1 = start of match) and then wraps the
# freshly-inserted text with bookmarks on both sides so it survives as its
# own run through any later edits.
function Insert-Shielded($searchText, $collapseDir, $newText, $shieldLeftName, $shieldRightName) {
    $r = Get-Range $searchText
    $r.Collapse($collapseDir)
    if ($collapseDir -eq 0) {
        $r.InsertAfter($newText)
    } else {
        $r.InsertBefore($newText)
    }
    $nr = Get-Range $newText
    $leftPt = $nr.Duplicate
    $leftPt.Collapse(1)
    if ($shieldLeftName) { Shield-Point $leftPt $shieldLeftName }
    $rightPt = $nr.Duplicate
    $rightPt.Collapse(0)
    if ($shieldRightName) { Shield-Point $rightPt $shieldRightName }
}

# --- Pre-shield the pre-existing run boundaries that must survive our edits ---
$rOnDesign = Get-Range " on design"
$pA = $rOnDesign.Duplicate; $pA.Collapse(1); Shield-Point $pA "PreShA"
$pB = $rOnDesign.Duplicate; $pB.Collapse(0); Shield-Point $pB "PreShB"

$rMeeting = Get-Range " (in a meeting setting or otherwise)"
$pC = $rMeeting.Duplicate; $pC.Collapse(1); Shield-Point $pC "PreShC"
$pD = $rMeeting.Duplicate; $pD.Collapse(0); Shield-Point $pD "PreShD"

# --- "...would scale badly to even larger projects." -> "...to larger projects." ---
$d.Content.Find.Execute("scale badly to even larger projects", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "scale badly to larger projects", 2) | Out-Null

# --- "...not having proper interfaces is that..." -> "...interfaces done earlier on is that..." ---
Insert-Shielded " is that people felt" 1 " done earlier on" "ShA_L" "ShA_R"

# --- "...running around especially me..." -> "...running around to find things ... design stage especially me..." ---
Insert-Shielded "especially me" 1 "to find things which should, really, have been agreed at the design stage " "ShB_L" "ShB_R"

# --- Move the _GoBack bookmark from "Our group " to right before "is not ideal" ---
if ($d.Bookmarks.Exists("_GoBack")) { $d.Bookmarks("_GoBack").Delete() }

# --- "...code. This is not ideal..." -> "...code. This strategy is not ideal..." (bookmark lands here) ---
Insert-Shielded "is not ideal" 1 "strategy " "ShC_L" "_GoBack"

# --- Clean up temporary shield bookmarks (keep the real _GoBack) ---
foreach ($n in @("PreShA", "PreShB", "PreShC", "PreShD", "ShA_L", "ShA_R", "ShB_L", "ShB_R", "ShC_L")) {
    if ($d.Bookmarks.Exists($n)) { $d.Bookmarks($n).Delete() }
}
